$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E..T (E=5 .. T=20) values per row, in order:
# E, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T
$data = @{
    2  = @(3, 1, 0.9142543333333334, 2.742763, 0.1175834869881751, 0.1175834869881751, 3, 1, 13.788265, 41.364795, 0.7604887294515349, 0.760488729451535, 12.60598102539834, 113.453829228585, 0.08942091662411834, 0.08942091662411834)
    3  = @(3, 1, 0.9142543333333334, 2.742763, 0.1175834869881751, 0.1175834869881751, 3, 1, 2.118986, 6.356958, 0.1168722076973129, 0.1168722076973129, 1.937292132772667, 17.435629194954, 0.01374224171305628, 0.01374224171305628)
    4  = @(3, 1, 0.9142543333333334, 2.742763, 0.1175834869881751, 0.1175834869881751, 3, 1, 1.925545, 5.776635, 0.1062030118040055, 0.1062030118040055, 1.760437860278333, 15.843940742505, 0.01248772045656128, 0.01248772045656128)
    5  = @(3, 1, 0.9142543333333334, 2.742763, 0.1175834869881751, 0.1175834869881751, 2, 0.6666666666666666, 0.2979986666666667, 0.893996, 0.0164360510471466, 0.0164360510471466, 0.2724465723275556, 2.452019150948, 0.001932608194439143, 0.001932608194439143)
    6  = @(3, 1, 4.826051, 14.478153, 0.6206849497708361, 0.620684949770836, 3, 1, 13.788265, 41.364795, 0.7604887294515349, 0.760488729451535, 66.542870091515, 598.8858308236349, 0.4720239088409129, 0.4720239088409129)
    7  = @(3, 1, 4.826051, 14.478153, 0.6206849497708361, 0.620684949770836, 3, 1, 2.118986, 6.356958, 0.1168722076973129, 0.1168722076973129, 10.226334504286, 92.03701053857399, 0.07254082036421337, 0.07254082036421337)
    8  = @(3, 1, 4.826051, 14.478153, 0.6206849497708361, 0.620684949770836, 3, 1, 1.925545, 5.776635, 0.1062030118040055, 0.1062030118040055, 9.292778372794999, 83.63500535515499, 0.06591861104708065, 0.06591861104708066)
    9  = @(3, 1, 4.826051, 14.478153, 0.6206849497708361, 0.620684949770836, 2, 0.6666666666666666, 0.2979986666666667, 0.893996, 0.0164360510471466, 0.0164360510471466, 1.438156763265333, 12.943410869388, 0.01020160951862908, 0.01020160951862908)
    10 = @(3, 1, 1.986145, 5.958435, 0.2554407961214246, 0.2554407961214246, 3, 1, 13.788265, 41.364795, 0.7604887294515349, 0.760488729451535, 27.385493588425, 246.469442295825, 0.1942598464924707, 0.1942598464924707)
    11 = @(3, 1, 1.986145, 5.958435, 0.2554407961214246, 0.2554407961214246, 3, 1, 2.118986, 6.356958, 0.1168722076973129, 0.1168722076973129, 4.20861344897, 37.87752104073, 0.02985392977867009, 0.02985392977867009)
    12 = @(3, 1, 1.986145, 5.958435, 0.2554407961214246, 0.2554407961214246, 3, 1, 1.925545, 5.776635, 0.1062030118040055, 0.1062030118040055, 3.824411574025, 34.419704166225, 0.02712858188570821, 0.02712858188570821)
    13 = @(3, 1, 1.986145, 5.958435, 0.2554407961214246, 0.2554407961214246, 2, 0.6666666666666666, 0.2979986666666667, 0.893996, 0.0164360510471466, 0.0164360510471466, 0.5918685618066667, 5.326817056259999, 0.004198437964575501, 0.004198437964575501)
    14 = @(1, 0.3333333333333333, 0.04891300000000001, 0.146739, 0.006290767119564404, 0.006290767119564403, 3, 1, 13.788265, 41.364795, 0.7604887294515349, 0.760488729451535, 0.6744254059450001, 6.069828653505001, 0.004784057494033025, 0.004784057494033025)
    15 = @(1, 0.3333333333333333, 0.04891300000000001, 0.146739, 0.006290767119564404, 0.006290767119564403, 3, 1, 2.118986, 6.356958, 0.1168722076973129, 0.1168722076973129, 0.103645962218, 0.9328136599620001, 0.0007352158413731578, 0.0007352158413731578)
    16 = @(1, 0.3333333333333333, 0.04891300000000001, 0.146739, 0.006290767119564404, 0.006290767119564403, 3, 1, 1.925545, 5.776635, 0.1062030118040055, 0.1062030118040055, 0.094184182585, 0.847657643265, 0.0006680984146553478, 0.0006680984146553479)
    17 = @(1, 0.3333333333333333, 0.04891300000000001, 0.146739, 0.006290767119564404, 0.006290767119564403, 2, 0.6666666666666666, 0.2979986666666667, 0.893996, 0.0164360510471466, 0.0164360510471466, 0.01457600878266667, 0.131184079044, 0.0001033953695028719, 0.0001033953695028719)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        # Columns start at E which is column index 5
        $col = 5 + $i
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
